# Apply edits to HAWAII_2018 sheet: rename header columns, fix casing of
# "de/del/la/los" -> "De/Del/La/Los" in a handful of place names, correct a
# floating point rounding value, and remove the trailing footnote rows
# (123-127) that are no longer part of the cleaned dataset.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header row renames (row 1) ---
$ws.Cells.Item(1, 1).Value = "mx_state"
$ws.Cells.Item(1, 2).Value = "mx_municipality"
$ws.Cells.Item(1, 3).Value = "n_matriculas"
$ws.Cells.Item(1, 4).Value = "pct_matriculas"

# --- Capitalization fixes of "de/del/la/los" -> "De/Del/La/Los" ---
$ws.Cells.Item(8, 1).Value = "Ciudad De México"
$ws.Cells.Item(19, 1).Value = "Estado De México"
$ws.Cells.Item(24, 2).Value = "San Antonio La Isla"
$ws.Cells.Item(25, 2).Value = "San Felipe Del Progreso"
$ws.Cells.Item(35, 2).Value = "Valle De Santiago"
$ws.Cells.Item(39, 2).Value = "Acapulco De Juárez"
$ws.Cells.Item(40, 2).Value = "Chilpancingo De Los Bravo"
$ws.Cells.Item(43, 2).Value = "Técpan De Galeana"
$ws.Cells.Item(48, 2).Value = "Autlán De Navarro"
$ws.Cells.Item(52, 2).Value = "La Manzanilla De La Paz"
$ws.Cells.Item(53, 2).Value = "Lagos De Moreno"
$ws.Cells.Item(72, 2).Value = "Ixtlán Del Río"
$ws.Cells.Item(75, 2).Value = "Santa María Del Oro"
$ws.Cells.Item(80, 2).Value = "Chalcatongo De Hidalgo"
$ws.Cells.Item(81, 2).Value = "Heroica Ciudad De Ejutla De Crespo"
$ws.Cells.Item(83, 2).Value = "Heroica Ciudad De Juchitán De Zaragoza"
$ws.Cells.Item(84, 2).Value = "Miahuatlán De Porfirio Díaz"
$ws.Cells.Item(85, 2).Value = "Oaxaca De Juárez"
$ws.Cells.Item(94, 2).Value = "Cuayuca De Andrade"
$ws.Cells.Item(96, 2).Value = "Los Reyes De Juárez"
$ws.Cells.Item(97, 2).Value = "Mazapiltepec De Juárez"
$ws.Cells.Item(101, 2).Value = "Cadereyta De Montes"

# --- Floating point rounding correction ---
$ws.Cells.Item(29, 4).Value = 0.09090909090909093

# --- Remove trailing footnote rows (rows 123-127), shrinking the sheet
#     dimension down to A1:D121 ---
$ws.Range("A123:A127").EntireRow.Delete()
